$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace all Fitness values (column C, rows 2-252) from 7573 to 7293
$ws.Range("C2:C252").Value = 7293
